# sprint_39.xlsx - "Test Summary" sheet restructure
#
# Day 2 -> Day 3 had two blank spacer rows (13:14) separating the blocks;
# Day 3 -> Day 4 had an (implicit) two-row gap bigger than the standard
# one-row gap used elsewhere. Both gaps are collapsed to match the
# standard spacing used between the rest of the day blocks, which shifts
# every block from Day 3 onward upward.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two explicit blank spacer rows between "Day 2" and "Day 3".
$ws.Rows("13:14").Delete()

# Remove the (now) two extra blank rows between "Day 3" and "Day 4" so the
# gap matches the rest of the blocks.
$ws.Rows("20:21").Delete()

# Scroll back to the top-left and move the active selection to F11 (matches
# the saved view state in the workbook).
$ws.Range("F11").Select()
